$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.860.15'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '2.620.99'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.50'
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.00'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').Value = '2.636.07'
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('E10').Value = '  +5.56%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.347'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('D14').Value = '3.077.64'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '60.789.95'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.72'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '2.626.99'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '355.68'
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.63'
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.96'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '0.0₃0848'
$ws.Range('E28').Value = '  -2.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.38'
$ws.Range('E29').Value = '  -3.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.47'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.58'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '152.29'
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.90'
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.00'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.20'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.872'
$ws.Range('E37').Value = '  +3.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.50'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.849'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.43'
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.77'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '294.46'
$ws.Range('E42').Value = '  -5.05%  '
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.625'
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0556'
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.995'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.94'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.95'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0235'
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').Value = '2.004.75'
$ws.Range('E51').Value = '  -2.41%  '
